$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for a new "DELETE /profile" endpoint block (3 content rows + 1 blank
#    separator) just above the existing "PUT /admin" block (currently at row 40).
#    Inserting 4 whole rows at row 40 shifts everything from row 40 downward by 4.
$ws.Rows("40:43").Insert()

# 2) The "DELETE /" (admin-delete-user) block used to be 2 rows (header + success).
#    It is being split into 3 rows (header + forbidden + success), so insert one more
#    row right after the "username" row of that block (now at row 50).
$ws.Rows("51").Insert()

# 3) Fill in the new "DELETE /profile" endpoint block (rows 40-42)
$ws.Range("B40").Value = "DELETE"
$ws.Range("C40").Value = "/profile"
$ws.Range("D40").Value = "token"
$ws.Range("E40").Value = 404
$ws.Range("F40").Value = "user not found"

$ws.Range("E41").Value = 403
$ws.Range("F41").Value = "not logged in"

$ws.Range("E42").Value = 200
$ws.Range("F42").Value = "user successfully deleted"

# 4) Update the "DELETE /" -> "DELETE /admin" block (now starting at row 49)
$ws.Range("C49").Value = "/admin"

# Row 50 ("username" row) used to report the final 200 success; now it reports the
# 403 "not logged in" case, and the 200 success moves to the newly inserted row 51.
$ws.Range("E50").Value = 403
$ws.Range("F50").Value = "not logged in"

$ws.Range("E51").Value = 200
$ws.Range("F51").Value = "user successfully deleted"
